# Reorganize and update repo:
#  - add a new "MB Endpoint" column (E) with Yes/No values for the trial rows
#  - add a bold, 16pt "section" style used on J17 (mirrors the existing
#    B29 16pt heading style, but bold) and bump row 17's height to match
#  - move the active selection to E22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "MB Endpoint" -------------------------------------------------
# Header (no special style, matches C1/D1 which also carry no explicit style)
$ws.Range("E1").Value = "MB Endpoint"

# Row 16 is the only "Yes" - write it first so the new shared-string order is
# MB Endpoint(48), Yes(49), No(50), matching how the rows were authored.
$ws.Range("E16").Value = "Yes"

$ws.Range("E2").Value = "No"
$ws.Range("E3").Value = "No"
$ws.Range("E4").Value = "No"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "No"
$ws.Range("E7").Value = "No"
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "No"
$ws.Range("E10").Value = "No"
$ws.Range("E11").Value = "No"
$ws.Range("E12").Value = "No"
$ws.Range("E13").Value = "No"
$ws.Range("E14").Value = "No"
$ws.Range("E15").Value = "No"
$ws.Range("E17").Value = "No"
$ws.Range("E18").Value = "No"
$ws.Range("E19").Value = "No"
$ws.Range("E20").Value = "No"
$ws.Range("E21").Value = "No"

# Column E (rows 2-21) picks up the same text style already used by columns
# A-C in this table (cellXf with fontId for 11pt black Aptos Narrow).
$ws.Range("A2").Copy()
$ws.Range("E2:E21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- J17: new bold 16pt heading-style cell ---------------------------------------
# Start from the existing 16pt/black/Aptos style (currently only used by B29)
# then make it bold - this mirrors an existing font/xf pair plus one new
# (bold) variant, instead of fabricating an unrelated style from scratch.
$ws.Range("B29").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J17").Font.Bold = $true

# Row 17 grows to fit the new 16pt content (matches row 29, which already
# hosts the same 16pt font and is 22pt tall).
$ws.Rows(17).RowHeight = 22

# --- Selection ---------------------------------------------------------------
$ws.Range("E22").Select() | Out-Null
